$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @{
    "B2" = 0.1826923076923077
    "C2" = 0.5673076923076923
    "J2" = 0.01682692307692308
    "P2" = 0.1418269230769231
    "S2" = 0.09134615384615384
    "B3" = 0.01666666666666667
    "C3" = 0.008333333333333333
    "J3" = 0.0125
    "P3" = 0.6916666666666667
    "S3" = 0.2708333333333333
    "J4" = 0.01886792452830189
    "O4" = 0.01886792452830189
    "P4" = 0.7358490566037735
    "S4" = 0.2264150943396226
    "B6" = 0.06072874493927125
    "D6" = 0.02024291497975709
    "F6" = 0.048582995951417
    "J6" = 0.2145748987854251
    "O6" = 0.01214574898785425
    "Q6" = 0.1336032388663968
    "R6" = 0.06882591093117409
    "S6" = 0.4412955465587045
    "B7" = 0.1355932203389831
    "D7" = 0.0211864406779661
    "F7" = 0.05084745762711865
    "J7" = 0.08050847457627118
    "O7" = 0.01694915254237288
    "Q7" = 0.1440677966101695
    "R7" = 0.05932203389830509
    "S7" = 0.4915254237288136
    "B8" = 0.1294326241134752
    "D8" = 0.01773049645390071
    "E8" = 0.003546099290780142
    "F8" = 0.06560283687943262
    "J8" = 0.09042553191489362
    "O8" = 0.01773049645390071
    "Q8" = 0.1436170212765958
    "R8" = 0.07624113475177305
    "S8" = 0.4556737588652482
    "B9" = 0.1582278481012658
    "D9" = 0.03164556962025317
    "F9" = 0.06329113924050633
    "J9" = 0.0949367088607595
    "O9" = 0.006329113924050633
    "Q9" = 0.1518987341772152
    "R9" = 0.05063291139240506
    "S9" = 0.4430379746835443
    "B10" = 0.1473354231974922
    "D10" = 0.02351097178683386
    "F10" = 0.08307210031347963
    "J10" = 0.103448275862069
    "O10" = 0.01410658307210031
    "Q10" = 0.1724137931034483
    "R10" = 0.06661442006269593
    "S10" = 0.3894984326018809
    "G11" = 0.1465968586387434
    "J11" = 0.1020942408376963
    "K11" = 0.2172774869109948
    "L11" = 0.5209424083769634
    "S11" = 0.01308900523560209
    "G12" = 0.719047619047619
    "J12" = 0.1761904761904762
    "K12" = 0.009523809523809525
    "L12" = 0.02380952380952381
    "S12" = 0.07142857142857142
    "G13" = 0.74
    "J13" = 0.24
    "S13" = 0.02
    "G14" = 0.5
    "S14" = 0.5
    "F15" = 0.03240740740740741
    "H15" = 0.1527777777777778
    "I15" = 0.07407407407407407
    "J15" = 0.3564814814814815
    "K15" = 0.06018518518518518
    "M15" = 0.009259259259259259
    "N15" = 0.004629629629629629
    "O15" = 0.05092592592592592
    "S15" = 0.2592592592592592
    "F16" = 0.01945525291828794
    "H16" = 0.1867704280155642
    "I16" = 0.07003891050583658
    "J16" = 0.3618677042801556
    "K16" = 0.1167315175097276
    "M16" = 0.01945525291828794
    "O16" = 0.05836575875486381
    "S16" = 0.1673151750972763
    "F17" = 0.01275510204081633
    "H17" = 0.2219387755102041
    "I17" = 0.06377551020408163
    "J17" = 0.3928571428571428
    "K17" = 0.09183673469387756
    "M17" = 0.02551020408163265
    "O17" = 0.03316326530612245
    "S17" = 0.1581632653061225
    "F18" = 0.03012048192771084
    "H18" = 0.2108433734939759
    "I18" = 0.04216867469879518
    "J18" = 0.4156626506024096
    "K18" = 0.1265060240963855
    "O18" = 0.0783132530120482
    "S18" = 0.0963855421686747
    "F19" = 0.01153846153846154
    "H19" = 0.2307692307692308
    "I19" = 0.05961538461538462
    "J19" = 0.3416666666666667
    "K19" = 0.1269230769230769
    "M19" = 0.02243589743589744
    "N19" = 0.000641025641025641
    "O19" = 0.06602564102564103
    "S19" = 0.1403846153846154
}

foreach ($cellRef in $changes.Keys) {
    $ws.Range($cellRef).Value = $changes[$cellRef]
}
